$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G2 (Base unit) remains "Box" - re-assign to keep it present in shared strings
$ws.Range("G2").Value = "Box"

# New columns: Selling Description, Purchase Description (added to shared strings
# right after "Box", matching the order new strings were introduced upstream)
$ws.Range("J2").Value = "Selling Desc"
$ws.Range("K2").Value = "Pur Desc"

# Update row 2 values (Name changes from "Kivi" to "Pear")
$ws.Range("A2").Value = "Pear"
$ws.Range("B2").Value = 400
$ws.Range("C2").Value = 401

# D2 (Brand) and E2 (Category) remain "General" / "Default" - unchanged

$ws.Range("F2").Value = 123

$ws.Range("H2").Value = 200
$ws.Range("I2").Value = 160

$ws.Range("L2").Value = 10
$ws.Range("M2").Value = 160
$ws.Range("N2").Value = 2

# Update sheet view: scroll so column B is the leftmost visible column,
# then select M2 to match the final selection state.
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("M2").Select()
